# Natmi following Dr Hou advice
# Update ligand/receptor-expressing cell counts and derived expression statistics
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 9.156959333333335
$ws.Range("H2").Value = 27.470878
$ws.Range("I2").Value = 0.969469463764299
$ws.Range("J2").Value = 0.9694694637642989
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.265744
$ws.Range("N2").Value = 0.7972319999999999
$ws.Range("O2").Value = 0.1034864391735229
$ws.Range("P2").Value = 0.1034864391735229
$ws.Range("Q2").Value = 2.433407001077334
$ws.Range("R2").Value = 21.900663009696
$ws.Range("S2").Value = 0.100326942692432
$ws.Range("T2").Value = 0.1003269426924319

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 9.156959333333335
$ws.Range("H3").Value = 27.470878
$ws.Range("I3").Value = 0.969469463764299
$ws.Range("J3").Value = 0.9694694637642989
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.307583
$ws.Range("N3").Value = 0.922749
$ws.Range("O3").Value = 0.1197794472135201
$ws.Range("P3").Value = 0.1197794472135201
$ws.Range("Q3").Value = 2.816525022624667
$ws.Range("R3").Value = 25.348725203622
$ws.Range("S3").Value = 0.1161225164600755
$ws.Range("T3").Value = 0.1161225164600755

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 9.156959333333335
$ws.Range("H4").Value = 27.470878
$ws.Range("I4").Value = 0.969469463764299
$ws.Range("J4").Value = 0.9694694637642989
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.994584333333333
$ws.Range("N4").Value = 5.983753
$ws.Range("O4").Value = 0.7767341136129571
$ws.Range("P4").Value = 0.7767341136129571
$ws.Range("Q4").Value = 18.26432762723712
$ws.Range("R4").Value = 164.378948645134
$ws.Range("S4").Value = 0.7530200046117916
$ws.Range("T4").Value = 0.7530200046117915

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.288371
$ws.Range("H5").Value = 0.865113
$ws.Range("I5").Value = 0.03053053623570109
$ws.Range("J5").Value = 0.03053053623570109
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.265744
$ws.Range("N5").Value = 0.7972319999999999
$ws.Range("O5").Value = 0.1034864391735229
$ws.Range("P5").Value = 0.1034864391735229
$ws.Range("Q5").Value = 0.07663286302399999
$ws.Range("R5").Value = 0.689695767216
$ws.Range("S5").Value = 0.003159496481090917
$ws.Range("T5").Value = 0.003159496481090916

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.288371
$ws.Range("H6").Value = 0.865113
$ws.Range("I6").Value = 0.03053053623570109
$ws.Range("J6").Value = 0.03053053623570109
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.307583
$ws.Range("N6").Value = 0.922749
$ws.Range("O6").Value = 0.1197794472135201
$ws.Range("P6").Value = 0.1197794472135201
$ws.Range("Q6").Value = 0.08869801729299999
$ws.Range("R6").Value = 0.798282155637
$ws.Range("S6").Value = 0.003656930753444621
$ws.Range("T6").Value = 0.003656930753444621

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.288371
$ws.Range("H7").Value = 0.865113
$ws.Range("I7").Value = 0.03053053623570109
$ws.Range("J7").Value = 0.03053053623570109
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.994584333333333
$ws.Range("N7").Value = 5.983753
$ws.Range("O7").Value = 0.7767341136129571
$ws.Range("P7").Value = 0.7767341136129571
$ws.Range("Q7").Value = 0.5751802787876666
$ws.Range("R7").Value = 5.176622509089
$ws.Range("S7").Value = 0.02371410900116555
$ws.Range("T7").Value = 0.02371410900116555

